$d = $word.ActiveDocument

# 1. Update the CV date line: "Sep 20, 2025" -> "Jan 17, 2026"
$d.Content.Find.Execute("Sep 20, 2025", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Jan 17, 2026", 2)

# 2. Update the personal site URL: "https://prajitdas.github.com" -> "https://prajitdas.github.io"
$d.Content.Find.Execute("prajitdas.github.com", $true, $false, $false, $false, $false,
                         $true, 1, $false, "prajitdas.github.io", 2)

# 3. Update the summary line to add "Security researcher and " before "Software Engineering Leader"
$d.Content.Find.Execute("As a Software Engineering Leader in Cisco", $true, $false, $false, $false, $false,
                         $true, 1, $false, "As a Security researcher and Software Engineering Leader in Cisco", 2)
